$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Target cluster labels
$ws.Range("D2").Value = "M1"
$ws.Range("D3").Value = "M2"

# Update row 2 values (D2 -> M1)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.9119773333333333
$ws.Range("H2").Value = 2.735932
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.8058999999999999
$ws.Range("N2").Value = 2.4177
$ws.Range("O2").Value = 0.09451456273674329
$ws.Range("P2").Value = 0.0945145627367433
$ws.Range("Q2").Value = 0.7349625329333332
$ws.Range("R2").Value = 6.6146627964
$ws.Range("S2").Value = 0.09451456273674329
$ws.Range("T2").Value = 0.0945145627367433

# Update row 3 values (D3 -> M2)
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.9119773333333333
$ws.Range("H3").Value = 2.735932
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.375759
$ws.Range("N3").Value = 7.127277
$ws.Range("O3").Value = 0.2786249200308754
$ws.Range("P3").Value = 0.2786249200308755
$ws.Range("Q3").Value = 2.166638357462666
$ws.Range("R3").Value = 19.499745217164
$ws.Range("S3").Value = 0.2786249200308754
$ws.Range("T3").Value = 0.2786249200308755

# Add new row 4 (D4 -> Neutro)
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Ccl20"
$ws.Range("C4").Value = "Cxcr3"
$ws.Range("D4").Value = "Neutro"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.9119773333333333
$ws.Range("H4").Value = 2.735932
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.345069333333334
$ws.Range("N4").Value = 16.035208
$ws.Range("O4").Value = 0.6268605172323812
$ws.Range("P4").Value = 0.6268605172323812
$ws.Range("Q4").Value = 4.874582077095112
$ws.Range("R4").Value = 43.871238693856
$ws.Range("S4").Value = 0.6268605172323812
$ws.Range("T4").Value = 0.6268605172323812
